# The workbook's single data table (rows 2..355) holds one price
# observation per row. This edit adds one new weekly observation: a new
# row is inserted right before the current row 250, which pushes the
# existing rows 250..355 down to 251..356, and the newly opened row 250
# is populated with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250 (shifts rows 250:355 down to 251:356).
$ws.Rows(250).Insert()

# Populate the new row with the new data point.
$ws.Range("A250").Value = 10
$ws.Range("B250").Value = "Vega Modelo de Temuco"
$ws.Range("C250").Value = "La Araucanía"
$ws.Range("D250").Value = 44704
$ws.Range("E250").Value = 9
$ws.Range("F250").Value = 100112037
$ws.Range("G250").Value = "Cebollín"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 50
$ws.Range("K250").Value = 9000
$ws.Range("L250").Value = 9000
$ws.Range("M250").Value = 9000
$ws.Range("N250").Value = "`$/docena de paquetes"
$ws.Range("O250").Value = "Provincia de Cautín"
$ws.Range("P250").Value = 750
$ws.Range("Q250").Value = 12
$ws.Range("R250").Value = "Hortaliza"
